$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the time values in column A (times were running ~38 min behind -> bump forward)
$ws.Range("A1").Value = 0.46736111111111112
$ws.Range("A2").Value = 0.4694444444444445
$ws.Range("A3").Value = 0.46666666666666662
$ws.Range("A4").Value = 0.46875

# Move the active selection from B4 to B3
$ws.Range("B3").Select()
